$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F6").Value = 2324
$ws.Range("F10").Value = 69
$ws.Range("F15").Value = 438
$ws.Range("F16").Value = 885
$ws.Range("F17").Value = 475
$ws.Range("F18").Value = 3221
$ws.Range("F19").Value = 406
$ws.Range("F20").Value = 130
$ws.Range("F21").Value = 3240
$ws.Range("F22").Value = 703
$ws.Range("F23").Value = 587
$ws.Range("F24").Value = 265
$ws.Range("F25").Value = 1059
$ws.Range("F28").Value = 872
$ws.Range("F29").Value = 842

# Sheet 2: 演出 (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F13").Value = 90
$ws.Range("F14").Value = 174
$ws.Range("G14").Value = 588
$ws.Range("F18").Value = 98
$ws.Range("F19").Value = 215
$ws.Range("F21").Value = 456

# Sheet 3: 本地生活 (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 240
$ws.Range("F6").Value = 448

# Sheet 4: 全部类型 (All types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 240
$ws.Range("F12").Value = 448
$ws.Range("F13").Value = 2324
$ws.Range("F21").Value = 69
$ws.Range("F30").Value = 438
$ws.Range("F31").Value = 90
$ws.Range("F32").Value = 885
$ws.Range("F33").Value = 475
$ws.Range("F34").Value = 3221
$ws.Range("F35").Value = 130
$ws.Range("F36").Value = 3240
$ws.Range("F37").Value = 703
$ws.Range("F39").Value = 587
$ws.Range("F40").Value = 265
$ws.Range("F41").Value = 1059
$ws.Range("F43").Value = 98
$ws.Range("F44").Value = 215
$ws.Range("F46").Value = 456
$ws.Range("F49").Value = 872
$ws.Range("F50").Value = 842
